$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Range("M28").Value = 1.1
$ws.Range("O28").Value = 1.5

# Row 29
$ws.Range("I29").Value = 3.7
$ws.Range("M29").Value = 1.04
$ws.Range("O29").Value = 1.22
$ws.Range("W29").Value = 9
$ws.Range("Z29").Value = 19
$ws.Range("AC29").Value = 13
$ws.Range("AE29").Value = 12
$ws.Range("AI29").Value = 21
$ws.Range("AJ29").Value = 13

# Row 35
$ws.Range("M35").Value = 1.01
$ws.Range("O35").Value = 1.11
$ws.Range("U35").Value = 1.67

# Row 36
$ws.Range("M36").Value = 1.03
$ws.Range("N36").Value = 13
$ws.Range("O36").Value = 1.22
$ws.Range("Q36").Value = 1.83
$ws.Range("R36").Value = 2.03
$ws.Range("U36").Value = 1.67

# Row 37
$ws.Range("M37").Value = 1.05
$ws.Range("O37").Value = 1.33
$ws.Range("U37").Value = 1.87
$ws.Range("V37").Value = 1.87

# Row 51
$ws.Range("Q51").Value = 2.05
$ws.Range("R51").Value = 1.75

# Row 75
$ws.Range("G75").Value = 2
$ws.Range("I75").Value = 3.8
$ws.Range("J75").Value = 2.75
$ws.Range("Y75").Value = 9
$ws.Range("Z75").Value = 17
$ws.Range("AI75").Value = 19
$ws.Range("BA75").Value = 81

# Row 96
$ws.Range("M96").Value = 1.06
$ws.Range("O96").Value = 1.3
$ws.Range("U96").Value = 1.8
$ws.Range("V96").Value = 1.91

# Row 97
$ws.Range("M97").Value = 1.03
$ws.Range("O97").Value = 1.2
$ws.Range("U97").Value = 1.67

# Row 102
$ws.Range("M102").Value = 1.06
$ws.Range("O102").Value = 1.3

# Row 104
$ws.Range("G104").Value = 1.57
$ws.Range("H104").Value = 3.9
$ws.Range("I104").Value = 4.55
$ws.Range("J104").Value = 2.07
$ws.Range("K104").Value = 2.32
$ws.Range("L104").Value = 4.6
$ws.Range("O104").Value = 1.14
$ws.Range("P104").Value = 4.15
$ws.Range("Q104").Value = 1.57
$ws.Range("R104").Value = 2.12
$ws.Range("U104").Value = 1.63
$ws.Range("V104").Value = 2.22
$ws.Range("X104").Value = 7.4
$ws.Range("Z104").Value = 10.25
$ws.Range("AA104").Value = 9.75
$ws.Range("AB104").Value = 17
$ws.Range("AC104").Value = 13.5
$ws.Range("AD104").Value = 6.9
$ws.Range("AH104").Value = 13
$ws.Range("AI104").Value = 23
$ws.Range("AJ104").Value = 12.5
$ws.Range("AK104").Value = 60
$ws.Range("AN104").Value = 3.6
$ws.Range("AO104").Value = 7.4
$ws.Range("AQ104").Value = 22
$ws.Range("AR104").Value = 45
$ws.Range("AT104").Value = 3.2
$ws.Range("AX104").Value = 6.5
$ws.Range("AY104").Value = 24
$ws.Range("BA104").Value = 120

# Row 105
$ws.Range("M105").Value = 46
$ws.Range("N105").Value = 26

# Row 106
$ws.Range("Q106").Value = 1.85
$ws.Range("R106").Value = 2
$ws.Range("AS106").Value = 151

# Row 110
$ws.Range("M110").Value = 1.03
$ws.Range("O110").Value = 1.17

# Row 111
$ws.Range("U111").Value = 1.5

# Row 112
$ws.Range("U112").Value = 1.95
$ws.Range("V112").Value = 1.8

# Row 118
$ws.Range("M118").Value = 1.06
$ws.Range("O118").Value = 1.3

# Row 119
$ws.Range("M119").Value = 1.05
$ws.Range("O119").Value = 1.25

# Row 120
$ws.Range("M120").Value = 1.05
$ws.Range("O120").Value = 1.25

# Row 121
$ws.Range("M121").Value = 1.05
$ws.Range("O121").Value = 1.25

# Row 122
$ws.Range("M122").Value = 1.03
$ws.Range("O122").Value = 1.18

# Row 123
$ws.Range("M123").Value = 1.03
$ws.Range("O123").Value = 1.2

# Row 136
$ws.Range("G136").Value = 2.05
$ws.Range("I136").Value = 3.5
$ws.Range("L136").Value = 3.75
$ws.Range("M136").Value = 1.05
$ws.Range("N136").Value = 11
$ws.Range("Q136").Value = 1.85
$ws.Range("R136").Value = 2
$ws.Range("X136").Value = 11
$ws.Range("Z136").Value = 19
$ws.Range("AD136").Value = 6.5
$ws.Range("AJ136").Value = 12
$ws.Range("AL136").Value = 26
$ws.Range("AN136").Value = 4.33
$ws.Range("AP136").Value = 21
$ws.Range("BB136").Value = 67

# Row 137
$ws.Range("M137").Value = 1.13
$ws.Range("N137").Value = 6

# Row 141
$ws.Range("Q141").Value = 1.98
$ws.Range("R141").Value = 1.88

# Row 142
$ws.Range("Q142").Value = 1.95
$ws.Range("R142").Value = 1.9

# Row 150
$ws.Range("M150").Value = 1.05
$ws.Range("N150").Value = 11
$ws.Range("O150").Value = 1.29
$ws.Range("P150").Value = 3.5
$ws.Range("Q150").Value = 2
$ws.Range("R150").Value = 1.85

# Row 158
$ws.Range("K158").Value = 1.91
$ws.Range("O158").Value = 1.53
$ws.Range("P158").Value = 2.38
$ws.Range("AL158").Value = 26

# Row 160
$ws.Range("O160").Value = 1.22
$ws.Range("P160").Value = 4
$ws.Range("Q160").Value = 1.73
$ws.Range("R160").Value = 2.08

# Row 161
$ws.Range("N161").Value = 9

Write-Output "Applied 130 cell updates across 26 rows."